# Generate Report for Handback
# Update the timestamp values recorded in the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 11:12:36"

# "zh-cn" sheet: Correspond Handoff / Handback Datetime for the first file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 11:12:31"
$wsZhCn.Range("K2").Value = "2016-08-28 11:12:48"

# "de-de" sheet: Correspond Handoff Datetime (shares the same underlying
# string as Overview!G2 in the source workbook, so it moves together with
# it) and Correspond Handback DateTime for the first file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 11:12:36"
$wsDeDe.Range("K2").Value = "2016-08-28 11:12:55"
